$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 312.75
$ws.Range("I9").Value = 351.3
$ws.Range("K9").Value = 351.3
$ws.Range("M9").Value = -182.3

$ws.Range("H43").Value = 3984.9285
$ws.Range("I43").Value = 2073
$ws.Range("J43").Value = 4506.364
$ws.Range("K43").Value = 2073
$ws.Range("L43").Value = 4506.364
$ws.Range("M43").Value = -2004
$ws.Range("N43").Value = -4644.364

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H49").Value = 878.3333
$ws.Range("I49").Value = 858
$ws.Range("J49").Value = 919
$ws.Range("K49").Value = 2574
$ws.Range("L49").Value = 2757
$ws.Range("M49").Value = -2438
$ws.Range("N49").Value = -3029

$ws.Range("H76").Value = 4470
$ws.Range("I76").Value = 4010
$ws.Range("K76").Value = 4010
$ws.Range("M76").Value = -3695

$ws.Range("H79").Value = 4470
$ws.Range("I79").Value = 4010
$ws.Range("K79").Value = 4010
$ws.Range("M79").Value = -2918

$ws.Range("H81").Value = 68146.664
$ws.Range("J81").Value = 68146.664
$ws.Range("L81").Value = 68146.664
$ws.Range("N81").Value = -70142.664

$ws.Range("H84").Value = 68146.664
$ws.Range("J84").Value = 68146.664
$ws.Range("L84").Value = 204439.992
$ws.Range("N84").Value = -214423.992

$ws.Range("H112").Value = 4812.324
$ws.Range("J112").Value = 5103.853
$ws.Range("L112").Value = 15311.559
$ws.Range("N112").Value = -17527.559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 877.4
$ws.Range("I63").Value = 971.75
$ws.Range("K63").Value = 971.75
$ws.Range("M63").Value = -285.75

$ws.Range("H66").Value = 877.4
$ws.Range("I66").Value = 971.75
$ws.Range("K66").Value = 4858.75
$ws.Range("M66").Value = -1426.75

$ws.Range("H97").Value = 2240.7917
$ws.Range("I97").Value = 1547.5555
$ws.Range("J97").Value = 2656.7334
$ws.Range("K97").Value = 1547.5555
$ws.Range("L97").Value = 2656.7334
$ws.Range("M97").Value = -1051.5555
$ws.Range("N97").Value = -3648.7334

$ws.Range("H110").Value = 437.5263
$ws.Range("I110").Value = 437
$ws.Range("J110").Value = 442
$ws.Range("K110").Value = 437
$ws.Range("L110").Value = 442
$ws.Range("M110").Value = 1608
$ws.Range("N110").Value = -4532

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1381.5172
$ws.Range("I80").Value = 1230.4667
$ws.Range("J80").Value = 1543.3572
$ws.Range("K80").Value = 1230.4667
$ws.Range("L80").Value = 1543.3572
$ws.Range("M80").Value = -232.4666999999999
$ws.Range("N80").Value = -3539.3572

$ws.Range("H83").Value = 1381.5172
$ws.Range("I83").Value = 1230.4667
$ws.Range("J83").Value = 1543.3572
$ws.Range("K83").Value = 6152.3335
$ws.Range("L83").Value = 7716.786
$ws.Range("M83").Value = -1160.3335
$ws.Range("N83").Value = -17700.786

$ws.Range("H86").Value = 3048.7222
$ws.Range("I86").Value = 1997.6923
$ws.Range("J86").Value = 5781.4
$ws.Range("K86").Value = 1997.6923
$ws.Range("L86").Value = 5781.4
$ws.Range("M86").Value = -874.6922999999999
$ws.Range("N86").Value = -8027.4

$ws.Range("H89").Value = 3048.7222
$ws.Range("I89").Value = 1997.6923
$ws.Range("J89").Value = 5781.4
$ws.Range("K89").Value = 9988.461499999999
$ws.Range("L89").Value = 28907
$ws.Range("M89").Value = -4372.461499999999
$ws.Range("N89").Value = -40139

$ws.Range("H94").Value = 1224.9166
$ws.Range("J94").Value = 1833
$ws.Range("L94").Value = 1833
$ws.Range("N94").Value = -2735

$ws.Range("H105").Value = 9705.583000000001
$ws.Range("I105").Value = 11215.6
$ws.Range("J105").Value = 2155.5
$ws.Range("K105").Value = 11215.6
$ws.Range("L105").Value = 2155.5
$ws.Range("M105").Value = -9468.6
$ws.Range("N105").Value = -5649.5

$ws.Range("H140").Value = 102044
$ws.Range("J140").Value = 102044
$ws.Range("L140").Value = 102044
$ws.Range("N140").Value = -112404

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8335883
$ws.Range("I31").Value = 12501564
$ws.Range("J31").Value = 4520.25
$ws.Range("K31").Value = 12501564
$ws.Range("L31").Value = 4520.25
$ws.Range("M31").Value = -12501269
$ws.Range("N31").Value = -5110.25

$ws.Range("H34").Value = 8335883
$ws.Range("I34").Value = 12501564
$ws.Range("J34").Value = 4520.25
$ws.Range("K34").Value = 12501564
$ws.Range("L34").Value = 4520.25
$ws.Range("M34").Value = -12501362
$ws.Range("N34").Value = -4924.25

$ws.Range("H58").Value = 20619
$ws.Range("I58").Value = 1604.8
$ws.Range("J58").Value = 47782.145
$ws.Range("K58").Value = 1604.8
$ws.Range("L58").Value = 47782.145
$ws.Range("M58").Value = -1401.8
$ws.Range("N58").Value = -48188.145

$ws.Range("H103").Value = 28749.75
$ws.Range("I103").Value = 17000
$ws.Range("K103").Value = 17000
$ws.Range("M103").Value = -15828

$ws.Range("H107").Value = 735.35297
$ws.Range("I107").Value = 435.8
$ws.Range("J107").Value = 1163.2858
$ws.Range("K107").Value = 435.8
$ws.Range("L107").Value = 1163.2858
$ws.Range("M107").Value = 1484.2
$ws.Range("N107").Value = -5003.2858

$ws.Range("H136").Value = 20619
$ws.Range("I136").Value = 1604.8
$ws.Range("J136").Value = 47782.145
$ws.Range("K136").Value = 4814.4
$ws.Range("L136").Value = 143346.435
$ws.Range("M136").Value = -2264.4
$ws.Range("N136").Value = -148446.435

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1389
$ws.Range("J92").Value = 443.5
$ws.Range("L92").Value = 1330.5
$ws.Range("N92").Value = -3826.5

$ws.Range("H138").Value = 7775.857
$ws.Range("J138").Value = 10144.6
$ws.Range("L138").Value = 30433.8
$ws.Range("N138").Value = -40713.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 118.21429
$ws.Range("I2").Value = 173.75
$ws.Range("J2").Value = 44.166668
$ws.Range("K2").Value = 173.75
$ws.Range("L2").Value = 44.166668
$ws.Range("M2").Value = -60.75
$ws.Range("N2").Value = -270.166668

$ws.Range("H58").Value = 9000
$ws.Range("J58").Value = 9000
$ws.Range("L58").Value = 9000
$ws.Range("N58").Value = -9554

$ws.Range("H62").Value = 85000
$ws.Range("J62").Value = 85000
$ws.Range("L62").Value = 85000
$ws.Range("N62").Value = -86372

$ws.Range("H65").Value = 85000
$ws.Range("J65").Value = 85000
$ws.Range("L65").Value = 255000
$ws.Range("N65").Value = -261864

$ws.Range("H126").Value = 2654.2727
$ws.Range("I126").Value = 2654.2727
$ws.Range("K126").Value = 7962.8181
$ws.Range("M126").Value = -5492.8181

$ws.Range("H132").Value = 4057.889
$ws.Range("I132").Value = 4324.7856
$ws.Range("K132").Value = 12974.3568
$ws.Range("M132").Value = -10444.3568

$ws.Range("H135").Value = 119997
$ws.Range("J135").Value = 119997
$ws.Range("L135").Value = 119997
$ws.Range("N135").Value = -130137

$ws.Range("H137").Value = 78001
$ws.Range("J137").Value = 78001
$ws.Range("L137").Value = 78001
$ws.Range("N137").Value = -88201

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3009.6
$ws.Range("I100").Value = 3093.25
$ws.Range("J100").Value = 2675
$ws.Range("K100").Value = 3093.25
$ws.Range("L100").Value = 2675
$ws.Range("M100").Value = -2552.25
$ws.Range("N100").Value = -3757

$ws.Range("H134").Value = 118844.336
$ws.Range("J134").Value = 118844.336
$ws.Range("L134").Value = 118844.336
$ws.Range("N134").Value = -128984.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 33752.168
$ws.Range("J105").Value = 33752.168
$ws.Range("L105").Value = 33752.168
$ws.Range("N105").Value = -40740.168

$ws.Range("H107").Value = 563.7895
$ws.Range("I107").Value = 474.2
$ws.Range("J107").Value = 899.75
$ws.Range("K107").Value = 1422.6
$ws.Range("L107").Value = 2699.25
$ws.Range("M107").Value = 497.4000000000001
$ws.Range("N107").Value = -6539.25

$ws.Range("H113").Value = 1412.7
$ws.Range("I113").Value = 1287.3334
$ws.Range("J113").Value = 1600.75
$ws.Range("K113").Value = 3862.0002
$ws.Range("L113").Value = 4802.25
$ws.Range("M113").Value = -1692.0002
$ws.Range("N113").Value = -9142.25

$ws.Range("H126").Value = 3727.3076
$ws.Range("I126").Value = 3316.3
$ws.Range("J126").Value = 5097.3335
$ws.Range("K126").Value = 9948.900000000001
$ws.Range("L126").Value = 15292.0005
$ws.Range("M126").Value = -7478.900000000001
$ws.Range("N126").Value = -20232.0005

$ws.Range("H133").Value = 82544
$ws.Range("J133").Value = 83555
$ws.Range("L133").Value = 83555
$ws.Range("N133").Value = -93675

$ws.Range("H136").Value = 11870.333
$ws.Range("I136").Value = 13620.102
$ws.Range("K136").Value = 40860.306
$ws.Range("M136").Value = -38310.306
